$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $row, $col, $val)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws 2 4 '41.591.25'
$ws.Cells.Item(2, 5).Value = '  +0.05%  '

Set-TextCell $ws 3 4 '2.470.05'
$ws.Cells.Item(3, 5).Value = '  +0.33%  '

Set-TextCell $ws 4 4 '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

Set-TextCell $ws 5 4 '318.65'
$ws.Cells.Item(5, 5).Value = '  +1.38%  '

Set-TextCell $ws 6 4 '92.34'
$ws.Cells.Item(6, 5).Value = '  +1.42%  '

Set-TextCell $ws 7 4 '0.552'
$ws.Cells.Item(7, 5).Value = '  +0.70%  '

$ws.Cells.Item(8, 5).Value = '  +0.05%  '

$ws.Cells.Item(9, 5).Value = '  +0.24%  '

Set-TextCell $ws 10 4 '0.0864'
$ws.Cells.Item(10, 5).Value = '  +8.80%  '

Set-TextCell $ws 11 4 '33.03'
$ws.Cells.Item(11, 5).Value = '  +1.87%  '

$ws.Cells.Item(12, 5).Value = '  -0.02%  '

Set-TextCell $ws 13 4 '2.850.88'
$ws.Cells.Item(13, 5).Value = '  +0.32%  '

Set-TextCell $ws 14 4 '6.87'
$ws.Cells.Item(14, 5).Value = '  +0.55%  '

Set-TextCell $ws 15 4 '15.47'
$ws.Cells.Item(15, 5).Value = '  -1.85%  '

Set-TextCell $ws 16 4 '2.471.27'
$ws.Cells.Item(16, 5).Value = '  +2.55%  '

Set-TextCell $ws 17 4 '0.791'
$ws.Cells.Item(17, 5).Value = '  +2.69%  '

Set-TextCell $ws 18 4 '41.563.51'
$ws.Cells.Item(18, 5).Value = '  +0.05%  '

Set-TextCell $ws 19 4 '6.44'
$ws.Cells.Item(19, 5).Value = '  -0.53%  '

$ws.Cells.Item(20, 5).Value = '  +0.70%  '

Set-TextCell $ws 21 4 '70.72'
$ws.Cells.Item(21, 5).Value = '  -0.31%  '

Set-TextCell $ws 22 4 '11.29'
$ws.Cells.Item(22, 5).Value = '  -0.13%  '

Set-TextCell $ws 23 4 '240.39'
$ws.Cells.Item(23, 5).Value = '  +1.46%  '

$ws.Cells.Item(24, 5).Value = '  +1.46%  '

$ws.Cells.Item(25, 5).Value = '  +3.41%  '

$ws.Cells.Item(26, 5).Value = '  +0.02%  '

Set-TextCell $ws 27 4 '24.73'
$ws.Cells.Item(27, 5).Value = '  +2.06%  '

Set-TextCell $ws 28 4 '2.24'
$ws.Cells.Item(28, 5).Value = '  -0.90%  '

$ws.Cells.Item(29, 5).Value = '  +0.59%  '

Set-TextCell $ws 30 4 '36.44'
$ws.Cells.Item(30, 5).Value = '  +4.09%  '

Set-TextCell $ws 31 4 '157.37'
$ws.Cells.Item(31, 5).Value = '  +0.93%  '

$ws.Cells.Item(32, 5).Value = '  +0.48%  '

$ws.Cells.Item(33, 5).Value = '  -0.01%  '

$ws.Cells.Item(34, 5).Value = '  +0.84%  '

$ws.Cells.Item(35, 5).Value = '  -0.69%  '

Set-TextCell $ws 36 4 '17.24'
$ws.Cells.Item(36, 5).Value = '  -0.13%  '

$ws.Cells.Item(37, 5).Value = '  +4.50%  '

$ws.Cells.Item(38, 2).Value = 'Stellar'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws 38 4 '0.116'
$ws.Cells.Item(38, 5).Value = '  +1.74%  '

$ws.Cells.Item(39, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws 39 4 '2.89'
$ws.Cells.Item(39, 5).Value = '  +1.09%  '

$ws.Cells.Item(40, 5).Value = '  +2.08%  '

$ws.Cells.Item(41, 2).Value = 'RenderToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws 41 4 '3.98'
$ws.Cells.Item(41, 5).Value = '  +0.69%  '

$ws.Cells.Item(42, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell $ws 42 4 '2.49'
$ws.Cells.Item(42, 5).Value = '  +3.76%  '

Set-TextCell $ws 43 4 '1.987.61'
$ws.Cells.Item(43, 5).Value = '  +1.60%  '

Set-TextCell $ws 45 4 '18.83'
$ws.Cells.Item(45, 5).Value = '  +0.69%  '

$ws.Cells.Item(46, 5).Value = '  +2.36%  '

Set-TextCell $ws 47 4 '9.45'
$ws.Cells.Item(47, 5).Value = '  +5.47%  '

Set-TextCell $ws 48 4 '2.707.65'
$ws.Cells.Item(48, 5).Value = '  +0.25%  '

Set-TextCell $ws 49 4 '97.46'
$ws.Cells.Item(49, 5).Value = '  +1.22%  '

Set-TextCell $ws 50 4 '75.91'
$ws.Cells.Item(50, 5).Value = '  +6.08%  '

Set-TextCell $ws 51 4 '66.88'
$ws.Cells.Item(51, 5).Value = '  +0.22%  '
